# EspecificacionLoginMiniCtrl.docx — "entrega: nuevo diseño, implementación
# del verificar de archivos"
#
# Changes applied:
#  1. FA2 table (3rd table, "Usuario-contraseña no es encontrado"), row 5.2,
#     "Descripción" cell: drop the stray "y" ("... usuario y la
#     contraseña ..." -> "... usuario la contraseña ...").
#  2. Same row, "Regla" cell (previously an empty list-styled paragraph):
#     add the new lockout rule text.
#  3. Page setup: stamp the section as explicit portrait orientation.

$d = $word.ActiveDocument

# --- locate the third table (the "FA2" flow table) -------------------------
$tbl = $d.Tables.Item(3)

# Row index 3 is "5.2" (1 = header, 2 = "5.1", 3 = "5.2", 4 = trailing blank row)
$row = $tbl.Rows.Item(3)

# --- 1. fix the description text in that row --------------------------------
$descCell = $row.Cells.Item(3)
[void]$descCell.Range.Find.Execute(
    "usuario y la contraseña",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "usuario la contraseña",
    2
)

# --- 2. add the missing lockout rule text in the "Regla" cell --------------
$reglaCell = $row.Cells.Item(4)
$reglaCell.Range.Text = "Después de 3 intentos fallidos consecutivos,el usuario se bloquea por 10 minutos."

# --- 3. make the section orientation explicit (portrait) -------------------
$d.PageSetup.Orientation = 0
